# ---------------------------------------------------------------------------
# Nexial mobile-commands update
#
# 1) Insert a new "mobile" command column into the hidden '#system' sheet.
#    This is modeled as a manual column shift (R..AH -> S..AI) rather than a
#    native Insert so that unrelated sheet metadata (col widths, etc.) is
#    left untouched, matching the author's edit.
# 2) Populate the freed-up column R with the 'mobile' command list.
# 3) Insert a "mobile" row into column A's category list (the `target`
#    list), shifting "number".."xml" down by one row and writing "mobile"
#    into the gap.
# 4) Fix up every defined name whose range moved because of the above
#    shifts, and add the new "mobile" defined name.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

$lastRow = 150

# --- Step 1: shift columns AH(34) .. R(18) right by one column (S..AI) -----
# NOTE: always (re)write the destination cell -- including blanking it out
# when the source is empty -- so that columns are truly *shifted* rather
# than leaving stale values behind in rows where the source column happened
# to have fewer entries than the destination column previously did.
for ($c = 34; $c -ge 18; $c--) {
    for ($r = 1; $r -le $lastRow; $r++) {
        $v = $ws.Cells.Item($r, $c).Value()
        if ($v -eq $null) {
            $v = ""
        }
        $ws.Cells.Item($r, $c + 1).Value = $v
    }
}

# --- Step 2: populate new column R (18) with the 'mobile' command list -----
$mobileCommands = @(
    "mobile",
    "assertElementPresent(locator)",
    "click(locator)",
    "clickUntilNotFound(locator,waitMs,max)",
    "closeApp()",
    "orientation(mode)",
    "screenshot(file,locator)",
    "scroll(locator,direction)",
    "slide(start,end)",
    "type(locator,text)",
    "use(profile)",
    "zoom(start1,end1,start2,end2)"
)
for ($i = 0; $i -lt $mobileCommands.Length; $i++) {
    $ws.Cells.Item($i + 1, 18).Value = $mobileCommands[$i]
}
# column R (18) previously held the (now relocated) "number" list, which was
# longer than the new "mobile" list -- blank out the leftover tail rows.
for ($r = $mobileCommands.Length + 1; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 18).Value = ""
}

# --- Step 3: shift column A rows 18..34 down to 19..35, insert "mobile" ----
for ($r = 34; $r -ge 18; $r--) {
    $v = $ws.Cells.Item($r, 1).Value()
    if ($v -eq $null) {
        $v = ""
    }
    $ws.Cells.Item($r + 1, 1).Value = $v
}
$ws.Cells.Item(18, 1).Value = "mobile"

# --- Step 4: fix up defined names that shifted, and add the new one --------
$wb.Names.Item("number").RefersTo       = "='#system'!`$S`$2:`$S`$16"
$wb.Names.Item("pdf").RefersTo          = "='#system'!`$T`$2:`$T`$21"
$wb.Names.Item("rdbms").RefersTo        = "='#system'!`$U`$2:`$U`$9"
$wb.Names.Item("redis").RefersTo        = "='#system'!`$V`$2:`$V`$10"
$wb.Names.Item("sms").RefersTo          = "='#system'!`$W`$2:`$W`$2"
$wb.Names.Item("sound").RefersTo        = "='#system'!`$X`$2:`$X`$5"
$wb.Names.Item("ssh").RefersTo          = "='#system'!`$Y`$2:`$Y`$9"
$wb.Names.Item("step").RefersTo         = "='#system'!`$Z`$2:`$Z`$4"
$wb.Names.Item("step.inTime").RefersTo  = "='#system'!`$AA`$2:`$AA`$4"
$wb.Names.Item("target").RefersTo       = "='#system'!`$A`$2:`$A`$35"
$wb.Names.Item("tn.5250").RefersTo      = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("web").RefersTo          = "='#system'!`$AB`$2:`$AB`$150"
$wb.Names.Item("webalert").RefersTo     = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("webcookie").RefersTo    = "='#system'!`$AD`$2:`$AD`$10"
$wb.Names.Item("ws").RefersTo           = "='#system'!`$AG`$2:`$AG`$17"
$wb.Names.Item("ws.async").RefersTo     = "='#system'!`$AH`$2:`$AH`$8"
$wb.Names.Item("xml").RefersTo          = "='#system'!`$AI`$2:`$AI`$27"
$wb.Names.Item("word").RefersTo         = "='#system'!`$AF`$2:`$AF`$9"
$wb.Names.Item("webmail").RefersTo      = "='#system'!`$AE`$2:`$AE`$4"
$wb.Names.Add("mobile", "='#system'!`$R`$2:`$R`$12")

Write-Host "mobile command column + target row inserted; defined names updated"
